# Added register page step def, feature file and page object
# -----------------------------------------------------------
# This script updates the DSAlgoTestData workbook:
#   1. Replaces the loginData test rows (2-4) with new negative-login
#      test data and clears/extends the trailing rows (5-10).
#   2. Adds a new "registerData" worksheet at the end of the workbook
#      with register-page test data, and makes it the active sheet
#      (mirrors pythonCode no longer being the selected tab).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. loginData sheet
# ---------------------------------------------------------------
$loginData = $wb.Worksheets.Item("loginData")

# Row 2: sonali / user / Invalid Username and Password
$loginData.Range("A2").Value = "sonali"
$loginData.Range("B2").Value = "user"
$loginData.Range("C2").Value = "Invalid Username and Password"
$loginData.Rows.Item(2).RowHeight = 13.8

# Row 3: username / Dsalgo@1 / Invalid Username and Password
$loginData.Range("A3").Value = "username"
$loginData.Range("B3").Value = "Dsalgo@1"
$loginData.Range("C3").Value = "Invalid Username and Password"
$loginData.Rows.Item(3).RowHeight = 13.8

# Row 4: N2324435 / sdetbatch / Invalid Username and Password
$loginData.Range("A4").Value = "N2324435"
$loginData.Range("B4").Value = "sdetbatch"
$loginData.Range("C4").Value = "Invalid Username and Password"
$loginData.Rows.Item(4).RowHeight = 13.8

# Row 5: cleared out (kept styled, but empty) and shrunk to match new row height
$loginData.Range("A5:C5").ClearContents()
$loginData.Rows.Item(5).RowHeight = 13.8

# Rows 6-8: new blank rows
$loginData.Rows.Item(6).RowHeight = 13.8
$loginData.Rows.Item(7).RowHeight = 13.8
$loginData.Rows.Item(8).RowHeight = 13.8

# Row 9: blank row but keeps the data-row style (copied from row 5)
$loginData.Range("A9:C9").Style = $loginData.Range("A5").Style
$loginData.Rows.Item(9).RowHeight = 13.8

# Row 10: new blank row
$loginData.Rows.Item(10).RowHeight = 13.8

# Update the sheet's selected cell
$loginData.Range("G18").Select()

# ---------------------------------------------------------------
# 2. New registerData sheet (added after the last existing sheet)
# ---------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$registerData = $wb.Worksheets.Add($null, $lastSheet)
$registerData.Name = "registerData"

$registerData.Range("A1").Value = "username"
$registerData.Range("B1").Value = "password"
$registerData.Range("C1").Value = "confirmpassword"
$registerData.Range("A1:C1").Font.Bold = $true

$registerData.Range("A2").Value = "Sdet147"
$registerData.Range("B2").Value = "demo1234"
$registerData.Range("C2").Value = "demo1234"

$registerData.Columns.Item(1).ColumnWidth = 8.67
$registerData.Columns.Item(2).ColumnWidth = 16.94
$registerData.Columns.Item(3).ColumnWidth = 31.68

# Select/activate the new sheet - this also clears tabSelected on pythonCode
$registerData.Range("C5").Select()
$registerData.Activate()

Write-Host "registerData sheet added and loginData updated"
